# Scheduled runner update: refresh market-board derived price/profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) on each Leve
# sheet with the latest pulled values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1600.4286
$ws.Range("I2").Value = 1440.8
$ws.Range("K2").Value = 1440.8
$ws.Range("M2").Value = -1327.8
$ws.Range("H6").Value = 2628.125
$ws.Range("I6").Value = 165
$ws.Range("J6").Value = 6733.3335
$ws.Range("K6").Value = 495
$ws.Range("L6").Value = 20200.0005
$ws.Range("M6").Value = -383
$ws.Range("N6").Value = -20424.0005
$ws.Range("H98").Value = 4691.091
$ws.Range("I98").Value = 3012.1765
$ws.Range("J98").Value = 10399.4
$ws.Range("K98").Value = 3012.1765
$ws.Range("L98").Value = 10399.4
$ws.Range("M98").Value = -1514.1765
$ws.Range("N98").Value = -13395.4
$ws.Range("H106").Value = 1913.3334
$ws.Range("I106").Value = 1554
$ws.Range("J106").Value = 2632
$ws.Range("K106").Value = 1554
$ws.Range("L106").Value = 2632
$ws.Range("M106").Value = -923
$ws.Range("N106").Value = -3894
$ws.Range("H107").Value = 2267.3845
$ws.Range("I107").Value = 438.85715
$ws.Range("J107").Value = 4400.6665
$ws.Range("K107").Value = 438.85715
$ws.Range("L107").Value = 4400.6665
$ws.Range("M107").Value = 1481.14285
$ws.Range("N107").Value = -8240.666499999999
$ws.Range("H116").Value = 50646.668
$ws.Range("I116").Value = 22250
$ws.Range("K116").Value = 22250
$ws.Range("M116").Value = -18808
$ws.Range("H122").Value = 4691.091
$ws.Range("I122").Value = 3012.1765
$ws.Range("J122").Value = 10399.4
$ws.Range("K122").Value = 9036.529500000001
$ws.Range("L122").Value = 31198.2
$ws.Range("M122").Value = -6586.529500000001
$ws.Range("N122").Value = -36098.2
$ws.Range("H138").Value = 2634.25
$ws.Range("I138").Value = 1685.4
$ws.Range("J138").Value = 3729.077
$ws.Range("K138").Value = 5056.200000000001
$ws.Range("L138").Value = 11187.231
$ws.Range("M138").Value = 83.79999999999927
$ws.Range("N138").Value = -21467.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5010
$ws.Range("I61").Value = 5010
$ws.Range("K61").Value = 5010
$ws.Range("M61").Value = -4798
$ws.Range("H102").Value = 2152.111
$ws.Range("I102").Value = 2152.111
$ws.Range("K102").Value = 2152.111
$ws.Range("M102").Value = -530.1109999999999
$ws.Range("H136").Value = 5010
$ws.Range("I136").Value = 5010
$ws.Range("K136").Value = 15030
$ws.Range("M136").Value = -12480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2366.0833
$ws.Range("I86").Value = 2073.1428
$ws.Range("J86").Value = 2776.2
$ws.Range("K86").Value = 2073.1428
$ws.Range("L86").Value = 2776.2
$ws.Range("M86").Value = -950.1428000000001
$ws.Range("N86").Value = -5022.2
$ws.Range("H89").Value = 2366.0833
$ws.Range("I89").Value = 2073.1428
$ws.Range("J89").Value = 2776.2
$ws.Range("K89").Value = 10365.714
$ws.Range("L89").Value = 13881
$ws.Range("M89").Value = -4749.714
$ws.Range("N89").Value = -25113
$ws.Range("H99").Value = 85309.25
$ws.Range("I99").Value = 2318.6667
$ws.Range("J99").Value = 334281
$ws.Range("K99").Value = 2318.6667
$ws.Range("L99").Value = 334281
$ws.Range("M99").Value = -820.6667000000002
$ws.Range("N99").Value = -337277

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2351.2156
$ws.Range("I31").Value = 1497.8788
$ws.Range("K31").Value = 1497.8788
$ws.Range("M31").Value = -1202.8788
$ws.Range("H34").Value = 2351.2156
$ws.Range("I34").Value = 1497.8788
$ws.Range("K34").Value = 1497.8788
$ws.Range("M34").Value = -1295.8788
$ws.Range("H86").Value = 47623508
$ws.Range("I86").Value = 66670230
$ws.Range("J86").Value = 6699.8335
$ws.Range("K86").Value = 66670230
$ws.Range("L86").Value = 6699.8335
$ws.Range("M86").Value = -66669107
$ws.Range("N86").Value = -8945.833500000001
$ws.Range("H89").Value = 47623508
$ws.Range("I89").Value = 66670230
$ws.Range("J89").Value = 6699.8335
$ws.Range("K89").Value = 333351150
$ws.Range("L89").Value = 33499.1675
$ws.Range("M89").Value = -333345534
$ws.Range("N89").Value = -44731.1675
$ws.Range("H105").Value = 849.5
$ws.Range("I105").Value = 850.1667
$ws.Range("J105").Value = 847.5
$ws.Range("K105").Value = 850.1667
$ws.Range("L105").Value = 847.5
$ws.Range("M105").Value = 896.8333
$ws.Range("N105").Value = -4341.5
$ws.Range("H132").Value = 1325.4615
$ws.Range("I132").Value = 833.8570999999999
$ws.Range("J132").Value = 1899
$ws.Range("K132").Value = 2501.5713
$ws.Range("L132").Value = 5697
$ws.Range("M132").Value = 28.42870000000039
$ws.Range("N132").Value = -10757

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 166750000
$ws.Range("J37").Value = 166750000
$ws.Range("L37").Value = 500250000
$ws.Range("N37").Value = -500250224
$ws.Range("H40").Value = 95.28570999999999
$ws.Range("I40").Value = 108.4
$ws.Range("J40").Value = 62.5
$ws.Range("K40").Value = 433.6
$ws.Range("L40").Value = 250
$ws.Range("M40").Value = -364.6
$ws.Range("N40").Value = -388
$ws.Range("H46").Value = 75084.2
$ws.Range("I46").Value = 167912.33
$ws.Range("J46").Value = 13198.777
$ws.Range("K46").Value = 503736.99
$ws.Range("L46").Value = 39596.331
$ws.Range("M46").Value = -503645.99
$ws.Range("N46").Value = -39778.331
$ws.Range("H51").Value = 349.66666
$ws.Range("I51").Value = 349.66666
$ws.Range("K51").Value = 1048.99998
$ws.Range("M51").Value = -588.9999800000001
$ws.Range("H56").Value = 6832.696
$ws.Range("I56").Value = 6832.696
$ws.Range("K56").Value = 6832.696
$ws.Range("M56").Value = -6302.696
$ws.Range("H92").Value = 987.8570999999999
$ws.Range("I92").Value = 930.8570999999999
$ws.Range("K92").Value = 2792.5713
$ws.Range("M92").Value = -1544.5713
$ws.Range("H107").Value = 941.3913
$ws.Range("J107").Value = 945.1905
$ws.Range("L107").Value = 2835.5715
$ws.Range("N107").Value = -6675.5715
$ws.Range("H114").Value = 15385754
$ws.Range("I114").Value = 25001062
$ws.Range("J114").Value = 1261.4
$ws.Range("K114").Value = 75003186
$ws.Range("L114").Value = 3784.2
$ws.Range("M114").Value = -74999932
$ws.Range("N114").Value = -10292.2
$ws.Range("H115").Value = 2998.75
$ws.Range("J115").Value = 2998.75
$ws.Range("L115").Value = 8996.25
$ws.Range("N115").Value = -11346.25
$ws.Range("H127").Value = 2750
$ws.Range("J127").Value = 2750
$ws.Range("L127").Value = 8250
$ws.Range("N127").Value = -18170
$ws.Range("H131").Value = 19322292
$ws.Range("J131").Value = 72221.69
$ws.Range("L131").Value = 216665.07
$ws.Range("N131").Value = -226745.07
$ws.Range("H138").Value = 10193.333
$ws.Range("I138").Value = 580
$ws.Range("J138").Value = 15000
$ws.Range("K138").Value = 1740
$ws.Range("L138").Value = 45000
$ws.Range("M138").Value = 3400
$ws.Range("N138").Value = -55280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1368
$ws.Range("I22").Value = 1035.625
$ws.Range("K22").Value = 1035.625
$ws.Range("M22").Value = -740.625
$ws.Range("H27").Value = 1368
$ws.Range("I27").Value = 1035.625
$ws.Range("K27").Value = 1035.625
$ws.Range("M27").Value = -928.625
$ws.Range("H128").Value = 89999.5
$ws.Range("J128").Value = 89999.5
$ws.Range("L128").Value = 89999.5
$ws.Range("N128").Value = -99959.5
$ws.Range("H132").Value = 3790.2727
$ws.Range("I132").Value = 2671.1428
$ws.Range("K132").Value = 8013.428400000001
$ws.Range("M132").Value = -5483.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18099.6
$ws.Range("J81").Value = 18099.6
$ws.Range("L81").Value = 36199.2
$ws.Range("N81").Value = -38321.2
$ws.Range("H84").Value = 18099.6
$ws.Range("J84").Value = 18099.6
$ws.Range("L84").Value = 180996
$ws.Range("N84").Value = -191604
$ws.Range("H132").Value = 7304.926
$ws.Range("I132").Value = 7783.5
$ws.Range("J132").Value = 5199.2
$ws.Range("K132").Value = 23350.5
$ws.Range("L132").Value = 15597.6
$ws.Range("M132").Value = -20820.5
$ws.Range("N132").Value = -20657.6
$ws.Range("H138").Value = 88418.664
$ws.Range("J138").Value = 88418.664
$ws.Range("L138").Value = 88418.664
$ws.Range("N138").Value = -98698.664
